# Auto-generated cell updates derived from the target diff.
# Maps SheetName -> RowNumber -> ColumnLetter -> new value ($null means clear the cell).
$updates = [ordered]@{
    "ALC" = [ordered]@{
        9 = @{ "H" = 439; "I" = 250.14285; "K" = 250.14285; "M" = -81.14285000000001 }
        11 = @{ "H" = 137.08333; "I" = 137.08333; "K" = 137.08333; "M" = 2.916670000000011 }
        43 = @{ "H" = 41669308; "I" = 1539.8; "K" = 1539.8; "M" = -1470.8 }
        69 = @{ "H" = 6125; "J" = 6125; "L" = 18375; "N" = -20123 }
        70 = @{ "H" = 52635892; "I" = 2749.375; "J" = 90914540; "K" = 8248.125; "L" = 272743620; "M" = -7978.125; "N" = -272744160 }
        72 = @{ "H" = 6125; "J" = 6125; "L" = 55125; "N" = -63861 }
        73 = @{ "H" = 52635892; "I" = 2749.375; "J" = 90914540; "K" = 8248.125; "L" = 272743620; "M" = -7312.125; "N" = -272745492 }
        75 = @{ "H" = 333371400; "J" = 333371400; "L" = 333371400; "N" = -333373272 }
        78 = @{ "H" = 333371400; "J" = 333371400; "L" = 1000114200; "N" = -1000123560 }
        100 = @{ "H" = 7750; "I" = 0; "J" = 7750; "K" = 0; "L" = $null; "M" = 7750; "N" = -8832 }
        111 = @{ "H" = 1117.1666; "I" = 2499; "K" = 7497; "M" = -4430 }
        114 = @{ "H" = 0; "J" = 0; "L" = $null; "N" = 0 }
        138 = @{ "H" = 2026675; "J" = 3277507.5; "L" = 9832522.5; "N" = -9842802.5 }
    }
    "ARM" = [ordered]@{
        43 = @{ "H" = 25091.666; "J" = 25091.666; "L" = 25091.666; "N" = -25717.666 }
        61 = @{ "H" = 381922.03; "I" = 3200.9355; "K" = 3200.9355; "M" = -2988.9355 }
        63 = @{ "H" = 7926.7334; "I" = 5000.125; "K" = 5000.125; "M" = -4314.125 }
        66 = @{ "H" = 7926.7334; "I" = 5000.125; "K" = 25000.625; "M" = -21568.625 }
        102 = @{ "H" = 7507; "I" = 6366; "J" = 7818.1816; "K" = 6366; "L" = 7818.1816; "M" = -4744; "N" = -11062.1816 }
        110 = @{ "H" = 3036.25; "I" = 715.8333; "K" = 715.8333; "M" = 1329.1667 }
        136 = @{ "H" = 381922.03; "I" = 3200.9355; "K" = 9602.806500000001; "M" = -7052.806500000001 }
    }
    "BSM" = [ordered]@{
        99 = @{ "H" = 4313.7144; "I" = 2715.5; "K" = 2715.5; "M" = -1217.5 }
        105 = @{ "H" = 2938.889; "I" = 2000; "K" = 2000; "M" = -253 }
    }
    "CRP" = [ordered]@{
        31 = @{ "H" = 5878.3687; "I" = 2360.5715; "J" = 7930.4165; "K" = 2360.5715; "L" = 7930.4165; "M" = -2065.5715; "N" = -8520.416499999999 }
        34 = @{ "H" = 5878.3687; "I" = 2360.5715; "J" = 7930.4165; "K" = 2360.5715; "L" = 7930.4165; "M" = -2158.5715; "N" = -8334.416499999999 }
        105 = @{ "H" = 2406.96; "I" = 2142.6667; "K" = 2142.6667; "M" = -395.6667000000002 }
    }
    "CUL" = [ordered]@{
        13 = @{ "H" = 516.5; "I" = 516.5; "K" = 1549.5; "M" = -1381.5 }
        44 = @{ "H" = 647; "J" = 2500; "L" = 7500; "N" = -8296 }
        113 = @{ "H" = 1483.3572; "I" = 3454; "J" = 1154.9166; "K" = 10362; "L" = 3464.7498; "M" = -8192; "N" = -7804.7498 }
        114 = @{ "H" = 1702.5; "I" = 1375.7142; "J" = 3990; "K" = 4127.142599999999; "L" = 11970; "M" = -873.1425999999992; "N" = -18478 }
        122 = @{ "H" = 1519.5625; "J" = 1684.75; "L" = 15162.75; "N" = -20062.75 }
        129 = @{ "H" = 2715.5; "I" = 890.5454999999999; "J" = 4946; "K" = 2671.6365; "L" = 14838; "M" = 2328.3635; "N" = -24838 }
        133 = @{ "H" = 9777.777; "I" = 2000; "K" = 6000; "M" = -940 }
        137 = @{ "H" = 4583.75; "J" = 39995; "L" = 119985; "N" = -130185 }
    }
    "LTW" = [ordered]@{
        2 = @{ "H" = 348665.8; "I" = 348665.8; "K" = 348665.8; "M" = -348553.8 }
        61 = @{ "H" = 4157.5; "I" = 2465.2307; "K" = 2465.2307; "M" = -2263.2307 }
        68 = @{ "H" = 6175.294; "I" = 2796.6667; "J" = 8018.1816; "K" = 2796.6667; "L" = 8018.1816; "M" = -2047.6667; "N" = -9516.1816 }
        71 = @{ "H" = 6175.294; "I" = 2796.6667; "J" = 8018.1816; "K" = 13983.3335; "L" = 40090.908; "M" = -10239.3335; "N" = -47578.908 }
        100 = @{ "H" = 58828384; "I" = 55557800; "K" = 55557800; "M" = -55557259 }
        113 = @{ "H" = 4157.5; "I" = 2465.2307; "K" = 2465.2307; "M" = -295.2307000000001 }
        127 = @{ "H" = 120000; "J" = 120000; "L" = 120000; "N" = -129920 }
    }
    "WVR" = [ordered]@{
        62 = @{ "H" = 38482000; "I" = 500000000; "J" = 22166.75; "K" = 500000000; "L" = 22166.75; "M" = -499999376; "N" = -23414.75 }
        65 = @{ "H" = 38482000; "I" = 500000000; "J" = 22166.75; "K" = 2500000000; "L" = 110833.75; "M" = -2499996880; "N" = -117073.75 }
        100 = @{ "H" = 1497.6666; "I" = 1059.875; "J" = 5000; "K" = 2119.75; "L" = 10000; "M" = -1578.75; "N" = -11082 }
    }
}
$wb = $excel.ActiveWorkbook

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cols = $rows[$rowNum]
        foreach ($colLetter in $cols.Keys) {
            $newVal = $cols[$colLetter]
            $cell = $ws.Range("$colLetter$rowNum")
            if ($null -eq $newVal) {
                $cell.ClearContents()
            } else {
                $cell.Value = $newVal
            }
        }
    }
}

Write-Host "Done applying $($updates.Keys.Count) sheet(s) of updates."
